$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a date/time number format to the "Giờ chiếu" (showtime) column header
# and its data cell. This mirrors Excel re-deriving new cellXf entries for
# D3 (header) and D4 (data) from their existing styles plus the new numFmt.
$ws.Range("D3:D4").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"

# Fill in the first data row with a sample showtime record.
$ws.Range("B4").Value = "Cún Cưng Đại Náo Nhà Hát"
$ws.Range("C4").Value = "Beta Bắc Giang"
$ws.Range("D4").Value = 45347.833333333336
$ws.Range("E4").Value = 40000

# Leave the selection where the user last clicked.
$ws.Range("F7").Select()
